# Weekly update: insert a new week block (4 rows) at the top of the
# "Agricola del Norte S.A. de Arica - Pina / Caramelo" price-history
# table (rows 233-236), pushing the rest of the history down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 233; everything below (old 233:260)
# shifts down to 237:264, carrying its existing data/formatting intact.
$ws.Rows("233:236").Insert()

# --- Row 233 : Calidad "Especial" --------------------------------------
$ws.Cells.Item(233,1).Value  = 1
$ws.Cells.Item(233,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(233,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(233,4).Value  = 44918
$ws.Cells.Item(233,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(233,5).Value  = 15
$ws.Cells.Item(233,6).Value  = "Fruta"
$ws.Cells.Item(233,7).Value  = 100108
$ws.Cells.Item(233,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(233,9).Value  = 100108005
$ws.Cells.Item(233,10).Value = "Piña"
$ws.Cells.Item(233,11).Value = "Caramelo"
$ws.Cells.Item(233,12).Value = "Especial"
$ws.Cells.Item(233,13).Value = 400
$ws.Cells.Item(233,14).Value = 20000
$ws.Cells.Item(233,15).Value = 22000
$ws.Cells.Item(233,16).Value = 21000
$ws.Cells.Item(233,17).Value = "$/caja 10 unidades"
$ws.Cells.Item(233,18).Value = "Ecuador"
$ws.Cells.Item(233,19).Value = 2100
$ws.Cells.Item(233,20).Value = 10

# --- Row 234 : Calidad "Primera" ----------------------------------------
$ws.Cells.Item(234,1).Value  = 1
$ws.Cells.Item(234,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(234,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(234,4).Value  = 44918
$ws.Cells.Item(234,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(234,5).Value  = 15
$ws.Cells.Item(234,6).Value  = "Fruta"
$ws.Cells.Item(234,7).Value  = 100108
$ws.Cells.Item(234,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(234,9).Value  = 100108005
$ws.Cells.Item(234,10).Value = "Piña"
$ws.Cells.Item(234,11).Value = "Caramelo"
$ws.Cells.Item(234,12).Value = "Primera"
$ws.Cells.Item(234,13).Value = 400
$ws.Cells.Item(234,14).Value = 20000
$ws.Cells.Item(234,15).Value = 22000
$ws.Cells.Item(234,16).Value = 21000
$ws.Cells.Item(234,17).Value = "$/caja 12 unidades"
$ws.Cells.Item(234,18).Value = "Ecuador"
$ws.Cells.Item(234,19).Value = 1750
$ws.Cells.Item(234,20).Value = 12

# --- Row 235 : Calidad "Segunda" -----------------------------------------
$ws.Cells.Item(235,1).Value  = 1
$ws.Cells.Item(235,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(235,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(235,4).Value  = 44918
$ws.Cells.Item(235,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(235,5).Value  = 15
$ws.Cells.Item(235,6).Value  = "Fruta"
$ws.Cells.Item(235,7).Value  = 100108
$ws.Cells.Item(235,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(235,9).Value  = 100108005
$ws.Cells.Item(235,10).Value = "Piña"
$ws.Cells.Item(235,11).Value = "Caramelo"
$ws.Cells.Item(235,12).Value = "Segunda"
$ws.Cells.Item(235,13).Value = 400
$ws.Cells.Item(235,14).Value = 20000
$ws.Cells.Item(235,15).Value = 22000
$ws.Cells.Item(235,16).Value = 21000
$ws.Cells.Item(235,17).Value = "$/caja 14 unidades"
$ws.Cells.Item(235,18).Value = "Ecuador"
$ws.Cells.Item(235,19).Value = 1500
$ws.Cells.Item(235,20).Value = 14

# --- Row 236 : Calidad "Tercera" -----------------------------------------
$ws.Cells.Item(236,1).Value  = 1
$ws.Cells.Item(236,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(236,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(236,4).Value  = 44918
$ws.Cells.Item(236,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(236,5).Value  = 15
$ws.Cells.Item(236,6).Value  = "Fruta"
$ws.Cells.Item(236,7).Value  = 100108
$ws.Cells.Item(236,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(236,9).Value  = 100108005
$ws.Cells.Item(236,10).Value = "Piña"
$ws.Cells.Item(236,11).Value = "Caramelo"
$ws.Cells.Item(236,12).Value = "Tercera"
$ws.Cells.Item(236,13).Value = 400
$ws.Cells.Item(236,14).Value = 20000
$ws.Cells.Item(236,15).Value = 22000
$ws.Cells.Item(236,16).Value = 21000
$ws.Cells.Item(236,17).Value = "$/caja 16 unidades"
$ws.Cells.Item(236,18).Value = "Ecuador"
$ws.Cells.Item(236,19).Value = 1312
$ws.Cells.Item(236,20).Value = 16
